# Applies:
#  1. Highlights the "Aanliggende Tiles fixen ..." bullet green and adds a new
#     bullet "Wincondition" (also green) right after it, in the same list.
#  2. Appends two new log paragraphs at the end of the document, after the
#     "... 19.30 TRIGGERS WERKEN ..." paragraph:
#       "17.30 start 18.30 verder"
#       "12.30 start 14.00 Winning condition en scene change"

$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# 1. "Aanliggende Tiles fixen ..." bullet: add green highlight, then insert a
#    new "Wincondition" bullet right after it (same list/number/style, green).
# ---------------------------------------------------------------------------
$tilesPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Aanliggende Tiles fixen*") {
        $tilesPara = $p
    }
}

$tilesXml = "<w:p $wNs>" +
              "<w:pPr>" +
                "<w:pStyle w:val='ListParagraph'/>" +
                "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr>" +
                "<w:rPr><w:highlight w:val='green'/><w:lang w:val='nl-NL'/></w:rPr>" +
              "</w:pPr>" +
              "<w:r>" +
                "<w:rPr><w:highlight w:val='green'/><w:lang w:val='nl-NL'/></w:rPr>" +
                "<w:t>Aanliggende Tiles fixen (wat is aanliggend?)</w:t>" +
              "</w:r>" +
            "</w:p>" +
            "<w:p $wNs>" +
              "<w:pPr>" +
                "<w:pStyle w:val='ListParagraph'/>" +
                "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr>" +
                "<w:rPr><w:highlight w:val='green'/><w:lang w:val='nl-NL'/></w:rPr>" +
              "</w:pPr>" +
              "<w:r>" +
                "<w:rPr><w:highlight w:val='green'/><w:lang w:val='nl-NL'/></w:rPr>" +
                "<w:t>Wincondition</w:t>" +
              "</w:r>" +
            "</w:p>"

$tilesPara.Range.InsertXML($tilesXml)

# ---------------------------------------------------------------------------
# 2. Append two new paragraphs after "... TRIGGERS WERKEN ..." (the last
#    paragraph in the document body).
# ---------------------------------------------------------------------------
$triggersPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*TRIGGERS WERKEN*") {
        $triggersPara = $p
    }
}

$endBefore = $d.Content.End
$triggersPara.Range.InsertParagraphAfter()

$newPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -eq $endBefore) {
        $newPara = $p
    }
}

$logXml = "<w:p $wNs>" +
            "<w:pPr><w:rPr><w:lang w:val='nl-NL'/></w:rPr></w:pPr>" +
            "<w:r><w:rPr><w:lang w:val='nl-NL'/></w:rPr><w:t>17.30 start</w:t></w:r>" +
            "<w:r><w:rPr><w:lang w:val='nl-NL'/></w:rPr><w:t xml:space='preserve'> 18.30 verder</w:t></w:r>" +
          "</w:p>" +
          "<w:p $wNs>" +
            "<w:pPr><w:rPr><w:lang w:val='nl-NL'/></w:rPr></w:pPr>" +
            "<w:r><w:rPr><w:lang w:val='nl-NL'/></w:rPr><w:t>12.30 start</w:t></w:r>" +
            "<w:r><w:rPr><w:lang w:val='nl-NL'/></w:rPr><w:t xml:space='preserve'> 14.00 Winning condition en scene change</w:t></w:r>" +
          "</w:p>"

$newPara.Range.InsertXML($logXml)
